$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column L (match formatting of the preceding header cell)
$ws.Range("L1").Value = "Poupança"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill formulas for rows 2..101: Poupança = SUM(J:K) - SUM(B:I), with IFERROR guard
for ($r = 2; $r -le 101; $r++) {
    $ws.Range("L$r").Formula = "=IFERROR(SUM(J$r`:K$r) - SUM(B$r`:I$r), 0)"
}

$ws.Range("L101").Select()
